$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 192, pushing existing rows 192-204 down to 193-205
$ws.Rows.Item(192).Insert()

# Populate the newly inserted row 192 with the new record
$ws.Cells.Item(192, 1).Value = 4
$ws.Cells.Item(192, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(192, 3).Value = "Los Lagos"
$ws.Cells.Item(192, 4).Value = 44610
$ws.Cells.Item(192, 5).Value = 10
$ws.Cells.Item(192, 6).Value = "Fruta"
$ws.Cells.Item(192, 7).Value = 100108
$ws.Cells.Item(192, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(192, 9).Value = 100108005
$ws.Cells.Item(192, 10).Value = "Piña"
$ws.Cells.Item(192, 11).Value = "Caramelo"
$ws.Cells.Item(192, 12).Value = "Tercera"
$ws.Cells.Item(192, 13).Value = 160
$ws.Cells.Item(192, 14).Value = 20000
$ws.Cells.Item(192, 15).Value = 20000
$ws.Cells.Item(192, 16).Value = 20000
$ws.Cells.Item(192, 17).Value = "`$/caja 16 unidades"
$ws.Cells.Item(192, 18).Value = "Ecuador"
$ws.Cells.Item(192, 19).Value = 1250
$ws.Cells.Item(192, 20).Value = 16
